$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (columns E..T)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.902492
$ws.Range("H2").Value = 3.804984
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.09215
$ws.Range("N2").Value = 0.1843
$ws.Range("O2").Value = 0.01023932159503481
$ws.Range("P2").Value = 0.01017879497979401
$ws.Range("Q2").Value = 0.1753146378
$ws.Range("R2").Value = 0.7012585512
$ws.Range("S2").Value = 0.01023932159503481
$ws.Range("T2").Value = 0.01017879497979401

# Update row 3 values (columns E..H, O..T; I-N unchanged)
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.902492
$ws.Range("H3").Value = 3.804984
$ws.Range("O3").Value = 0.01189268776136058
$ws.Range("P3").Value = 0.01773358166721151
$ws.Range("Q3").Value = 0.203623084596
$ws.Range("R3").Value = 1.221738507576
$ws.Range("S3").Value = 0.01189268776136058
$ws.Range("T3").Value = 0.01773358166721151

# Update row 4 values (columns E..H, M..T; I-L unchanged)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.902492
$ws.Range("H4").Value = 3.804984
$ws.Range("M4").Value = 8.80044
$ws.Range("N4").Value = 17.60088
$ws.Range("O4").Value = 0.9778679906436047
$ws.Range("P4").Value = 0.9720876233529945
$ws.Range("Q4").Value = 16.74276669648
$ws.Range("R4").Value = 66.97106678592
$ws.Range("S4").Value = 0.9778679906436047
$ws.Range("T4").Value = 0.9720876233529945

# Remove rows 5 and 6 (Neutrophils / Resolving-Mac target rows no longer present)
$ws.Rows.Item(6).EntireRow.Delete() | Out-Null
$ws.Rows.Item(5).EntireRow.Delete() | Out-Null
